$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All values are assigned with a leading apostrophe to force Excel to
# store them as literal text (matching the workbook's inline-string cells)
# instead of auto-converting number-looking strings like '585.30' into 585.3.

$ws.Range("D2").Value = "'63.917.49"
$ws.Range("E2").Value = "'  -2.73%  "
$ws.Range("D3").Value = "'3.504.47"
$ws.Range("E3").Value = "'  -2.10%  "
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("D5").Value = "'585.30"
$ws.Range("E5").Value = "'  -3.12%  "
$ws.Range("D6").Value = "'131.93"
$ws.Range("E6").Value = "'  -3.73%  "
$ws.Range("D7").Value = "'3.503.68"
$ws.Range("E7").Value = "'  -2.09%  "
$ws.Range("E8").Value = "'  +0.07%  "
$ws.Range("E9").Value = "'  -2.10%  "
$ws.Range("D10").Value = "'0.124"
$ws.Range("E10").Value = "'  -0.96%  "
$ws.Range("D11").Value = "'7.12"
$ws.Range("E11").Value = "'  -1.20%  "
$ws.Range("D12").Value = "'0.386"
$ws.Range("E12").Value = "'  -1.35%  "
$ws.Range("D13").Value = "'4.107.67"
$ws.Range("E13").Value = "'  -2.04%  "
$ws.Range("D14").Value = "'27.77"
$ws.Range("E14").Value = "'  -1.42%  "
$ws.Range("D15").Value = "'0.0000180"
$ws.Range("E15").Value = "'  -3.11%  "
$ws.Range("E16").Value = "'  +0.76%  "
$ws.Range("D17").Value = "'3.514.26"
$ws.Range("D18").Value = "'63.998.85"
$ws.Range("E18").Value = "'  -2.69%  "
$ws.Range("D19").Value = "'10.00"
$ws.Range("E19").Value = "'  -0.65%  "
$ws.Range("D20").Value = "'14.48"
$ws.Range("E20").Value = "'  -1.14%  "
$ws.Range("D21").Value = "'5.66"
$ws.Range("E21").Value = "'  -3.50%  "
$ws.Range("D22").Value = "'390.95"
$ws.Range("E22").Value = "'  -0.85%  "
$ws.Range("D23").Value = "'0.579"
$ws.Range("E23").Value = "'  -1.72%  "
$ws.Range("D24").Value = "'3.650.51"
$ws.Range("E24").Value = "'  -1.97%  "
$ws.Range("D25").Value = "'73.08"
$ws.Range("E25").Value = "'  -1.34%  "
$ws.Range("E26").Value = "'  +0.01%  "
$ws.Range("E27").Value = "'  -3.77%  "
$ws.Range("E28").Value = "'  -0.37%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "'  +0.13%  "
$ws.Range("D31").Value = "'2.26"
$ws.Range("E31").Value = "'  -3.71%  "
$ws.Range("D32").Value = "'8.28"
$ws.Range("E32").Value = "'  -3.57%  "
$ws.Range("D33").Value = "'3.510.56"
$ws.Range("E33").Value = "'  -2.05%  "
$ws.Range("E34").Value = "'  +0.00%  "
$ws.Range("D35").Value = "'23.95"
$ws.Range("E35").Value = "'  -2.14%  "
$ws.Range("E36").Value = "'  -2.50%  "
$ws.Range("D37").Value = "'5.34"
$ws.Range("E37").Value = "'  -0.09%  "
$ws.Range("D38").Value = "'7.00"
$ws.Range("E38").Value = "'  -0.50%  "
$ws.Range("E39").Value = "'  -3.29%  "
$ws.Range("D40").Value = "'168.09"
$ws.Range("E40").Value = "'  +0.40%  "
$ws.Range("D41").Value = "'0.0811"
$ws.Range("E41").Value = "'  -2.73%  "
$ws.Range("D42").Value = "'26.93"
$ws.Range("E42").Value = "'  -0.26%  "
$ws.Range("D43").Value = "'0.813"
$ws.Range("E43").Value = "'  -2.94%  "
$ws.Range("E44").Value = "'  +0.12%  "
$ws.Range("D45").Value = "'41.93"
$ws.Range("E45").Value = "'  -2.66%  "

# Rows 46 and 47 swapped places (Filecoin moved above ONDO) with updated values
$ws.Range("B46").Value = "'Filecoin"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "'4.39"
$ws.Range("E46").Value = "'  -3.08%  "
$ws.Range("B47").Value = "'ONDO"
$ws.Range("C47").Value = "'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "'1.20"
$ws.Range("E47").Value = "'  -6.19%  "

$ws.Range("D48").Value = "'1.64"
$ws.Range("E48").Value = "'  -3.63%  "
$ws.Range("D49").Value = "'2.447.73"
$ws.Range("E49").Value = "'  -0.32%  "
$ws.Range("D51").Value = "'0.898"
$ws.Range("E51").Value = "'  -0.09%  "
